$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric ("319.10", "5.130", "0.09040", ...).
# Force it to Text format *before* writing so Excel keeps the exact string (incl.
# trailing zeros / precision) instead of silently coercing to a Double and losing
# formatting (e.g. "319.10" -> 319.1, "0.000009070" -> 9.07E-06).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.251.24'
$ws.Range("E2").Value = '  -2.45%  '
$ws.Range("D3").Value = '1.867.44'
$ws.Range("E3").Value = '  -1.90%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '319.10'
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = '0.4376'
$ws.Range("E7").Value = '  -4.64%  '
$ws.Range("D8").Value = '0.3698'
$ws.Range("E8").Value = '  -3.23%  '
$ws.Range("D9").Value = '0.07507'
$ws.Range("E9").Value = '  -2.59%  '
$ws.Range("D10").Value = '0.9381'
$ws.Range("E10").Value = '  -3.94%  '
$ws.Range("D11").Value = '21.44'
$ws.Range("E11").Value = '  -2.81%  '
$ws.Range("D12").Value = '1.866.93'
$ws.Range("E12").Value = '  -4.27%  '
$ws.Range("D13").Value = '6.718'
$ws.Range("E13").Value = '  -3.15%  '
$ws.Range("E14").Value = '  -3.77%  '
$ws.Range("D15").Value = '0.06874'
$ws.Range("E15").Value = '  -2.22%  '
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '82.07'
$ws.Range("E17").Value = '  -1.92%  '
$ws.Range("D18").Value = '0.000009070'
$ws.Range("E18").Value = '  -4.09%  '
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").Value = '15.92'
$ws.Range("E20").Value = '  -4.29%  '
$ws.Range("D21").Value = '28.242.52'
$ws.Range("E21").Value = '  -2.48%  '
$ws.Range("D22").Value = '5.130'
$ws.Range("E22").Value = '  -3.17%  '
$ws.Range("E23").Value = '  -0.86%  '
$ws.Range("D24").Value = '2.151.13'
$ws.Range("E24").Value = '  -0.85%  '
$ws.Range("D25").Value = '2.029'
$ws.Range("E25").Value = '  -3.21%  '
$ws.Range("D26").Value = '154.75'
$ws.Range("E26").Value = '  -2.14%  '
$ws.Range("D27").Value = '18.41'
$ws.Range("E27").Value = '  -3.28%  '
$ws.Range("D28").Value = '5.311'
$ws.Range("E28").Value = '  -5.96%  '
$ws.Range("D29").Value = '113.91'
$ws.Range("E29").Value = '  -2.93%  '
$ws.Range("D30").Value = '1.731'
$ws.Range("E30").Value = '  -5.87%  '
$ws.Range("D31").Value = '0.09040'
$ws.Range("E31").Value = '  -2.31%  '
$ws.Range("D32").Value = '0.7990'
$ws.Range("E32").Value = '  -7.52%  '
$ws.Range("D33").Value = '4.837'
$ws.Range("E33").Value = '  -4.92%  '
$ws.Range("D34").Value = '1.171'
$ws.Range("E34").Value = '  -5.67%  '
$ws.Range("D35").Value = '2.960'
$ws.Range("E35").Value = '  -1.24%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").Value = '1.122'
$ws.Range("E37").Value = '  -2.21%  '
$ws.Range("D38").Value = '0.05428'
$ws.Range("E38").Value = '  -5.00%  '
$ws.Range("D39").Value = '0.01952'
$ws.Range("E39").Value = '  -4.27%  '
$ws.Range("D40").Value = '2.929'
$ws.Range("E40").Value = '  +5.79%  '
$ws.Range("D41").Value = '7.105'
$ws.Range("E41").Value = '  -3.83%  '
$ws.Range("D42").Value = '0.5254'
$ws.Range("E42").Value = '  -4.32%  '
$ws.Range("D43").Value = '0.1676'
$ws.Range("E43").Value = '  -4.48%  '
$ws.Range("D44").Value = '8.721'
$ws.Range("E44").Value = '  -5.90%  '
$ws.Range("D45").Value = '0.06750'
$ws.Range("E45").Value = '  -0.97%  '
$ws.Range("D46").Value = '0.4875'
$ws.Range("E46").Value = '  -5.60%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '1.997'
$ws.Range("E47").Value = '  -3.04%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '10.57'
$ws.Range("E48").Value = '  -6.14%  '
$ws.Range("D49").Value = '107.83'
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("B50").Value = 'PEPE'
$ws.Range("C50").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D50").Value = '0.000002439'
$ws.Range("E50").Value = '  -6.05%  '
$ws.Range("D51").Value = '1.677'
$ws.Range("E51").Value = '  -5.34%  '
